$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.343283414840698
$ws.Range("B1").Value = 1.943539261817932
$ws.Range("C1").Value = 2.894254922866821
$ws.Range("D1").Value = 3.814670324325562
$ws.Range("E1").Value = 1.024856805801392
